$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '28.304.41'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.45%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.906.71'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +2.14%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '314.94'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5083'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.67%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3942'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.10%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.09679'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('E10').Value = '  +1.02%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '42.12'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +2.31%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '6.437'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.927.38'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '20.95'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.21%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.346'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.60%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.01%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.00001124'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.88%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '92.98'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.57%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06637'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.06%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '18.07'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +3.82%  '
$ws.Range('E21').Value = '  +0.11%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.229'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.00%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '28.360.91'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.40%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.34'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.14%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.316'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.00%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.671'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +4.73%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.136.27'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +2.62%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '21.02'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.02%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '157.98'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -0.03%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '127.58'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.46%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.098'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +4.25%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.1066'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.65%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.657'
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.634'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.60%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '9.710'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +3.01%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.06686'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.73%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.02434'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.249'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.2210'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.274'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +8.72%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.6404'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.031'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '11.53'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.76%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.23%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '13.49'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.00%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.6027'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.760'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.63%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.281'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +3.40%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '123.91'
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.193'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.53%  '
